$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.667632579803467
$ws.Range("B1").Value = 3.2921302318573
$ws.Range("C1").Value = 2.923240661621094
$ws.Range("D1").Value = 2.529649257659912
$ws.Range("E1").Value = 1.682997345924377
